$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9: new hours entry
$ws.Range("A9").Value = Get-Date -Year 2016 -Month 1 -Day 13 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("B9").Value = 1.5
$ws.Range("C9").Value = "Portaal tabs + nieuws op index en portaal"
